$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 232 (new weekly price entry, date 45006)
$ws.Rows("232:235").Insert()

# Row 232
$ws.Cells.Item(232, 1).Value = 6
$ws.Cells.Item(232, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(232, 3).Value = 'Metropolitana'
$ws.Cells.Item(232, 4).Value = 45006
$ws.Cells.Item(232, 5).Value = 13
$ws.Cells.Item(232, 6).Value = 100112043
$ws.Cells.Item(232, 7).Value = 'Pepino dulce'
$ws.Cells.Item(232, 8).Value = 'Cultivar IV Región'
$ws.Cells.Item(232, 9).Value = 'Especial'
$ws.Cells.Item(232, 10).Value = 280
$ws.Cells.Item(232, 11).Value = 15000
$ws.Cells.Item(232, 12).Value = 15000
$ws.Cells.Item(232, 13).Value = 15000
$ws.Cells.Item(232, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(232, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(232, 16).Value = 833
$ws.Cells.Item(232, 17).Value = 18
$ws.Cells.Item(232, 18).Value = 'Hortaliza'

# Row 233
$ws.Cells.Item(233, 1).Value = 6
$ws.Cells.Item(233, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(233, 3).Value = 'Metropolitana'
$ws.Cells.Item(233, 4).Value = 45006
$ws.Cells.Item(233, 5).Value = 13
$ws.Cells.Item(233, 6).Value = 100112043
$ws.Cells.Item(233, 7).Value = 'Pepino dulce'
$ws.Cells.Item(233, 8).Value = 'Cultivar IV Región'
$ws.Cells.Item(233, 9).Value = 'Primera'
$ws.Cells.Item(233, 10).Value = 470
$ws.Cells.Item(233, 11).Value = 13000
$ws.Cells.Item(233, 12).Value = 13000
$ws.Cells.Item(233, 13).Value = 13000
$ws.Cells.Item(233, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(233, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(233, 16).Value = 722
$ws.Cells.Item(233, 17).Value = 18
$ws.Cells.Item(233, 18).Value = 'Hortaliza'

# Row 234
$ws.Cells.Item(234, 1).Value = 6
$ws.Cells.Item(234, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(234, 3).Value = 'Metropolitana'
$ws.Cells.Item(234, 4).Value = 45006
$ws.Cells.Item(234, 5).Value = 13
$ws.Cells.Item(234, 6).Value = 100112043
$ws.Cells.Item(234, 7).Value = 'Pepino dulce'
$ws.Cells.Item(234, 8).Value = 'Cultivar IV Región'
$ws.Cells.Item(234, 9).Value = 'Segunda'
$ws.Cells.Item(234, 10).Value = 320
$ws.Cells.Item(234, 11).Value = 10000
$ws.Cells.Item(234, 12).Value = 10000
$ws.Cells.Item(234, 13).Value = 10000
$ws.Cells.Item(234, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(234, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(234, 16).Value = 556
$ws.Cells.Item(234, 17).Value = 18
$ws.Cells.Item(234, 18).Value = 'Hortaliza'

# Row 235
$ws.Cells.Item(235, 1).Value = 6
$ws.Cells.Item(235, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(235, 3).Value = 'Metropolitana'
$ws.Cells.Item(235, 4).Value = 45006
$ws.Cells.Item(235, 5).Value = 13
$ws.Cells.Item(235, 6).Value = 100112043
$ws.Cells.Item(235, 7).Value = 'Pepino dulce'
$ws.Cells.Item(235, 8).Value = 'Cultivar IV Región'
$ws.Cells.Item(235, 9).Value = 'Tercera'
$ws.Cells.Item(235, 10).Value = 180
$ws.Cells.Item(235, 11).Value = 8000
$ws.Cells.Item(235, 12).Value = 8000
$ws.Cells.Item(235, 13).Value = 8000
$ws.Cells.Item(235, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(235, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(235, 16).Value = 444
$ws.Cells.Item(235, 17).Value = 18
$ws.Cells.Item(235, 18).Value = 'Hortaliza'
